# The presentation's Handout Master and Notes Master each contain an
# automatically-updating "date" placeholder (a datetimeFigureOut field).
# The commit bumps the cached/displayed date shown by those fields from
# 5/15/2023 to 5/18/2023. Update both via the standard PowerPoint object
# model: Presentation.HandoutMaster / Presentation.NotesMaster expose the
# placeholder like any other shape's TextFrame.TextRange, and the
# Master's HeadersFooters.DateAndTime also mirrors the same value.

$p = $ppt.ActivePresentation

$newDate = "5/18/2023"

# --- Handout Master date placeholder ---
$handoutMaster = $p.HandoutMaster
for ($i = 1; $i -le $handoutMaster.Shapes.Count; $i++) {
    $shape = $handoutMaster.Shapes.Item($i)
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        if ($shape.TextFrame.TextRange.Text -eq "5/15/2023") {
            $shape.TextFrame.TextRange.Text = $newDate
        }
    }
}
$handoutMaster.HeadersFooters.DateAndTime.Value = $newDate

# --- Notes Master date placeholder ---
$notesMaster = $p.NotesMaster
for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
    $shape = $notesMaster.Shapes.Item($i)
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        if ($shape.TextFrame.TextRange.Text -eq "5/15/2023") {
            $shape.TextFrame.TextRange.Text = $newDate
        }
    }
}
$notesMaster.HeadersFooters.DateAndTime.Value = $newDate
